$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update BIT emission value for level 23 (row 24), which cascades through
# the dependent formulas in columns F, G, H, I, J for this and subsequent rows.
$ws.Range("E24").Value = 25800

# Restore the active selection on the sheet as left by the author.
$ws.Activate()
$ws.Range("I16").Select() | Out-Null
